# The document's headers carry the BTec logo picture (currently named
# "image1.jpg") and its footers carry the Pearson logo picture (currently
# named "image2.png"). The edit simply swaps those two picture names:
#   header logos: image1.jpg -> image2.jpg
#   footer logos: image2.png -> image1.png
#
# InlineShape has no writable "Name" in the Word object model (only a
# floating Shape does), so each picture is briefly converted to a Shape,
# renamed, and converted back to an InlineShape so the drawing stays
# inline exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

$headers = $sec.Headers
for ($h = 1; $h -le $headers.Count; $h++) {
    $hdr = $headers.Item($h)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -ge 1) {
        $inlineShp = $hdr.Range.InlineShapes.Item(1)
        $floatingShp = $inlineShp.ConvertToShape()
        $floatingShp.Name = "image2.jpg"
        $floatingShp.ConvertToInlineShape() | Out-Null
    }
}

$footers = $sec.Footers
for ($f = 1; $f -le $footers.Count; $f++) {
    $ftr = $footers.Item($f)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -ge 1) {
        $inlineShp = $ftr.Range.InlineShapes.Item(1)
        $floatingShp = $inlineShp.ConvertToShape()
        $floatingShp.Name = "image1.png"
        $floatingShp.ConvertToInlineShape() | Out-Null
    }
}

Write-Output "Renamed header (BTec) and footer (Pearson) logo images"
